$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remap F17:I17 off the soon-to-be-repurposed "style 8" onto the default style (they look identical - no fill).
$ws.Range("F17:I17").Interior.ColorIndex = 0

# New section: "#pragma pack(push, 1)" header (mirrors B13/K13 "Struct:" style row)
$ws.Range("B27").Value = "#pragma pack(push, 1)"
$ws.Range("K27").Value = "#pragma pack(push, 1)"

# Word-size sub-headers
$ws.Range("B29").Value = "64-bit word size"
$ws.Range("K29").Value = "32-bit word size"

$ws.Range("B29").Font.Bold = $true
$ws.Range("K29").Font.Bold = $true

# 64-bit grid (B31:I34), bytes 0-31, colored by field like the earlier grid (rows 15-18)
$ws.Range("B31").Value = 0
$ws.Range("C31").Value = 1
$ws.Range("D31").Value = 2
$ws.Range("E31").Value = 3
$ws.Range("F31").Value = 4
$ws.Range("G31").Value = 5
$ws.Range("H31").Value = 6
$ws.Range("I31").Value = 7

$ws.Range("B32").Value = 8
$ws.Range("C32").Value = 9
$ws.Range("D32").Value = 10
$ws.Range("E32").Value = 11
$ws.Range("F32").Value = 12
$ws.Range("G32").Value = 13
$ws.Range("H32").Value = 14
$ws.Range("I32").Value = 15

$ws.Range("B33").Value = 16
$ws.Range("C33").Value = 17
$ws.Range("D33").Value = 18
$ws.Range("E33").Value = 19
$ws.Range("F33").Value = 20
$ws.Range("G33").Value = 21
$ws.Range("H33").Value = 22
$ws.Range("I33").Value = 23

$ws.Range("B34").Value = 24
$ws.Range("C34").Value = 25
$ws.Range("D34").Value = 26
$ws.Range("E34").Value = 27
$ws.Range("F34").Value = 28
$ws.Range("G34").Value = 29
$ws.Range("H34").Value = 30
$ws.Range("I34").Value = 31

# 32-bit grid (K31:N37), bytes 0-27
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 1
$ws.Range("M31").Value = 2
$ws.Range("N31").Value = 3

$ws.Range("K32").Value = 4
$ws.Range("L32").Value = 5
$ws.Range("M32").Value = 6
$ws.Range("N32").Value = 7

$ws.Range("K33").Value = 8
$ws.Range("L33").Value = 9
$ws.Range("M33").Value = 10
$ws.Range("N33").Value = 11

$ws.Range("K34").Value = 12
$ws.Range("L34").Value = 13
$ws.Range("M34").Value = 14
$ws.Range("N34").Value = 15

$ws.Range("K35").Value = 16
$ws.Range("L35").Value = 17
$ws.Range("M35").Value = 18
$ws.Range("N35").Value = 19

$ws.Range("K36").Value = 20
$ws.Range("L36").Value = 21
$ws.Range("M36").Value = 22
$ws.Range("N36").Value = 23

$ws.Range("K37").Value = 24
$ws.Range("L37").Value = 25
$ws.Range("M37").Value = 26
$ws.Range("N37").Value = 27

# Color the byte grids to mirror the earlier layout pattern (struct field color bands).
$ws.Range("B31:E31").Interior.ColorIndex = 3
$ws.Range("F31").Interior.ColorIndex = 4
$ws.Range("G31:H31").Interior.ColorIndex = 5
$ws.Range("I31").Interior.ColorIndex = 6

$ws.Range("B32:D32").Interior.ColorIndex = 6
$ws.Range("E32:I32").Interior.ColorIndex = 7

$ws.Range("B33").Interior.ColorIndex = 7
$ws.Range("C33:E33").Interior.ColorIndex = 8
$ws.Range("F33:I33").Interior.ColorIndex = 9

$ws.Range("B34").Interior.ColorIndex = 8
$ws.Range("C34:I34").Interior.ColorIndex = 10

$ws.Range("K31:N31").Interior.ColorIndex = 3

$ws.Range("K32").Interior.ColorIndex = 4
$ws.Range("L32:M32").Interior.ColorIndex = 5
$ws.Range("N32").Interior.ColorIndex = 6

$ws.Range("K33:M33").Interior.ColorIndex = 6
$ws.Range("N33").Interior.ColorIndex = 7

$ws.Range("K34:N34").Interior.ColorIndex = 7

$ws.Range("K35").Interior.ColorIndex = 7
$ws.Range("L35:N35").Interior.ColorIndex = 8

$ws.Range("K36:N36").Interior.ColorIndex = 8

$ws.Range("K37").Interior.ColorIndex = 8
$ws.Range("L37:N37").Interior.ColorIndex = 10

# Summary rows
$ws.Range("B39").Value = "Size: 25 bytes"
$ws.Range("K39").Value = "Size: 25 bytes"
$ws.Range("B40").Value = "Holes: 0"
$ws.Range("K40").Value = "Holes: 0"

$ws.Range("K39").Font.Bold = $true
$ws.Range("K40").Font.Bold = $true

$ws.Range("B42").Select()
